$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 16: section header ---
$ws.Range("A16").Value = "EMBER Utilization Factors"

# --- Row 17: header row with fuel label + years 2000-2023 ---
$ws.Range("A17").Value = "model_fuel"
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(17, 2 + $i).Value = 2000 + $i
}

# --- Row 18: bioenergy ---
$ws.Range("A18").Value = "bioenergy"
$row18 = @(0, 0, 0, 0, 0, 0, 0, 0, 0.22831050228310504, 0.11415525114155252, 0.4566210045662101, 0.684931506849315, 0.7990867579908677, 0.4185692541856926, 0.5707762557077626, 0.6164383561643836, 0.6659056316590564, 0.9132420091324202, 2.5603392041748205, 3.4627092846270933, 3.881278538812785, 5.91324200913242, 5.114155251141553, 4.9771689497716896)
for ($i = 0; $i -lt $row18.Count; $i++) {
    $ws.Cells.Item(18, 2 + $i).Value = $row18[$i]
}

# --- Row 19: coal ---
$ws.Range("A19").Value = "coal"
$row19 = @(0.34409073920603195, 0.39609028421001313, 0.3489656965501552, 0.3908090804205463, 0.3922789539227895, 0.38314653383146535, 0.396222498962225, 0.46430053964300544, 0.48111249481112495, 0.42782873873654675, 0.4584458664849915, 0.49726171897578175, 0.4130902837986244, 0.36227010141320837, 0.41231328844516674, 0.5082841518219721, 0.4369935206742831, 0.46734400271649285, 0.41685655309224456, 0.384240767051801, 0.3018077187715018, 0.3817834133090278, 0.4867794368638804, 0.2582455387859779)
for ($i = 0; $i -lt $row19.Count; $i++) {
    $ws.Cells.Item(19, 2 + $i).Value = $row19[$i]
}

# --- Row 20: gas ---
$ws.Range("A20").Value = "gas"
$row20 = @(0.19642930601834707, 0.19642930601834707, 0.1583775556378296, 0.1810029207289481, 0.15323542720802993, 0.1779176436710683, 0.22213994816734542, 0.23228111971411552, 0.2342664284296208, 0.19456025411951558, 0.19555290847726822, 0.19786910197869106, 0.2245053272450533, 0.2226027397260274, 0.20357686453576868, 0.17694063926940642, 0.1950152207001522, 0.18359969558599698, 0.19216133942161343, 0.2045281582952816, 0.21784627092846273, 0.290144596651446, 0.1950152207001522, 0.14840182648401828)
for ($i = 0; $i -lt $row20.Count; $i++) {
    $ws.Cells.Item(20, 2 + $i).Value = $row20[$i]
}

# --- Row 21: hydro ---
$ws.Range("A21").Value = "hydro"
$row21 = @(0.2943414808845913, 0.22423352902804958, 0.22408253001860307, 0.20686315206863154, 0.1810340851436742, 0.24791291914579586, 0.2427240440939071, 0.1607260501147232, 0.1502326182476092, 0.18296846327828276, 0.2633949143311968, 0.14626141552511415, 0.15647142182333496, 0.19708855325293684, 0.22250599798777185, 0.27329541057193724, 0.1876789722157728, 0.12870890865760704, 0.2342229256490022, 0.13325692663137406, 0.1282541068602303, 0.219214466335571, 0.17145847997545438, 0.14032522966412186)
for ($i = 0; $i -lt $row21.Count; $i++) {
    $ws.Cells.Item(21, 2 + $i).Value = $row21[$i]
}

# --- Row 22: nuclear ---
$ws.Range("A22").Value = "nuclear"
$row22 = @(0.587915712678024, 0.632219592016247, 0.8486099919419822, 0.7252215954875101, 0.7059159280150417, 0.7827189094816007, 0.8179727370400214, 0.8842501993186925, 0.9525017515885096, 0.9216979536614239, 0.921093957623638, 0.985117537628953, 0.943125582729685, 0.8169595498362622, 0.9149716341497164, 0.8867210922005443, 0.9144009456922327, 0.9010731752543866, 0.9160816919966379, 0.9405029646289103, 0.9444785206388151, 0.9365274086190054, 0.934823598900475, 0.9177855017151686)
for ($i = 0; $i -lt $row22.Count; $i++) {
    $ws.Cells.Item(22, 2 + $i).Value = $row22[$i]
}

# --- Row 23: oil ---
$ws.Range("A23").Value = "oil"
$row23 = @(65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535, 65535)
for ($i = 0; $i -lt $row23.Count; $i++) {
    $ws.Cells.Item(23, 2 + $i).Value = $row23[$i]
}

# --- Row 24: solar ---
$ws.Range("A24").Value = "solar"
$row24 = @(0.03805175038051751, 0.07610350076103502, 0.09678379988088147, 0.1525728837372673, 0.1396462295518021, 0.15294587046149752, 0.1540541738706388, 0.1551624772797801, 0.1485126568249324, 0.15586582367404284, 0.15255292652552924, 0.132132456045734, 0.13711751430220961, 0.1366756748361445)
for ($i = 0; $i -lt $row24.Count; $i++) {
    $ws.Cells.Item(24, 12 + $i).Value = $row24[$i]
}

# --- Row 25: wind ---
$ws.Range("A25").Value = "wind"
$row25 = @(0, 0.07610350076103502, 0.19025875190258754, 0.12453300124533001, 0.08302200083022, 0.15841953219644023, 0.18180280737358362, 0.20480795057749127, 0.22998925597636316, 0.2168949771689498, 0.23646444879321593, 0.23157208088714937, 0.24461839530332682, 0.2152641878669276, 0.2152641878669276, 0.24135681669928247, 0.23320287018917157, 0.23483365949119375, 0.2527723418134377)
for ($i = 0; $i -lt $row25.Count; $i++) {
    $ws.Cells.Item(25, 7 + $i).Value = $row25[$i]
}

